$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.894.36'
$ws.Range("E2").Value = '  +4.14%  '
$ws.Range("D3").Value = '2.281.28'
$ws.Range("E3").Value = '  +4.71%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.74'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("E6").Value = '  +4.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.89'
$ws.Range("E7").Value = '  +9.73%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.654'
$ws.Range("E9").Value = '  +13.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.79'
$ws.Range("E10").Value = '  +6.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0977'
$ws.Range("E11").Value = '  +4.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.67'
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("E13").Value = '  +7.56%  '
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").Value = '2.623.81'
$ws.Range("E15").Value = '  +4.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.02'
$ws.Range("E16").Value = '  +4.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.889'
$ws.Range("E17").Value = '  +4.65%  '
$ws.Range("D18").Value = '2.283.56'
$ws.Range("E18").Value = '  +4.70%  '
$ws.Range("D19").Value = '42.834.17'
$ws.Range("E19").Value = '  +4.17%  '
$ws.Range("E20").Value = '  +6.96%  '
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.48'
$ws.Range("E22").Value = '  +2.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.21'
$ws.Range("E23").Value = '  +2.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").Value = '  +5.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.89'
$ws.Range("E25").Value = '  +2.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.68'
$ws.Range("E26").Value = '  +1.64%  '
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("E30").Value = '  +4.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.08'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.07'
$ws.Range("E32").Value = '  +4.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.36'
$ws.Range("E33").Value = '  +10.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.129'
$ws.Range("E34").Value = '  +5.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0813'
$ws.Range("E35").Value = '  +7.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '30.83'
$ws.Range("E36").Value = '  +25.61%  '
$ws.Range("E37").Value = '  +4.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.75'
$ws.Range("E38").Value = '  +20.63%  '
$ws.Range("E39").Value = '  +5.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0310'
$ws.Range("E40").Value = '  +1.46%  '
$ws.Range("E41").Value = '  +5.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.30'
$ws.Range("E42").Value = '  +16.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.98'
$ws.Range("E43").Value = '  +8.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.211'
$ws.Range("E44").Value = '  +11.74%  '
$ws.Range("E45").Value = '  +7.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.97'
$ws.Range("E46").Value = '  -6.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '61.39'
$ws.Range("E47").Value = '  +1.13%  '
$ws.Range("E48").Value = '  +3.56%  '
$ws.Range("E49").Value = '  +3.65%  '
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("E51").Value = '  +4.69%  '
